# Regenerate merged AHB files
# 1) Rename header labels from "_old"/"_new" suffixes to "_FV2410"/"_FV2504"
# 2) Freeze the header row (pane split) and keep the selection anchored below it
# 3) Convert the data range into a formatted Excel Table ("Table1")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update header row labels ---------------------------------------
$headers = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410",
    "diff",
    "Segmentname_FV2504",
    "Segmentgruppe_FV2504",
    "Segment_FV2504",
    "Datenelement_FV2504",
    "Segment ID_FV2504",
    "Code_FV2504",
    "Qualifier_FV2504",
    "Beschreibung_FV2504",
    "Bedingungsausdruck_FV2504",
    "Bedingung_FV2504"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2) Freeze panes at row 1 --------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true

# --- 3) Turn the used range into an Excel Table --------------------------
$dataRange = $ws.UsedRange
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"

$wb.Save()
